{"js": "// Author's edit: right after \"{{ fecha_dia }}\" and before \"{{ fecha_mes }}\",\n// the lone \" de\" run gets expanded into \" del mes n\u00famero \" so the sentence\n//   Medell\u00edn, {{ fecha_dia }} de{{ fecha_mes }} de {{ fecha_a\u00f1o }}\n// becomes\n//   Medell\u00edn, {{ fecha_dia }} del mes n\u00famero {{ fecha_mes }} de {{ fecha_a\u00f1o }}\n//\n// Locate the exact paragraph that contains the unedited pattern, then\n// replace just the \" de\" immediately preceding \"{{ fecha_mes }}\" (search\n// results come back in document order, so it is always the first match).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items.find(\n  (p) => p.text.indexOf(\"{{ fecha_dia }} de{{ fecha_mes }}\") !== -1\n);\n\nif (!targetParagraph) {\n  throw new Error(\"Could not find the Medell\u00edn date paragraph to edit.\");\n}\n\nconst matches = targetParagraph.search(\" de\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \" de\" inside the target paragraph.');\n}\n\n// First match = left-most occurrence = \" de\" right before \"{{ fecha_mes }}\".\nconst target = matches.items[0];\ntarget.insertText(\" del mes n\u00famero \", \"Replace\");\nawait context.sync();\n", "ps1": "# Author's edit: right after \"{{ fecha_dia }}\" and before \"{{ fecha_mes }}\",\n# the lone \" de\" text gets expanded into \" del mes n\u00famero \" so the sentence\n#   Medellin, {{ fecha_dia }} de{{ fecha_mes }} de {{ fecha_anio }}\n# becomes\n#   Medellin, {{ fecha_dia }} del mes n\u00famero {{ fecha_mes }} de {{ fecha_anio }}\n#\n# Locate the paragraph that still has the un-edited pattern, then replace just\n# the \" de\" that immediately precedes \"{{ fecha_mes }}\" (not the other \" de\"\n# later in the same paragraph, right before \"{{ fecha_a\u00f1o }}\").\n\n$d = $word.ActiveDocument\n\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"{{ fecha_dia }} de{{ fecha_mes }}\")) {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -eq $null) {\n    throw \"Could not find the Medellin date paragraph to edit.\"\n}\n\n$rng = $targetParagraph.Range\n$find = $rng.Find\n$find.Text = \" de\"\n$find.Replacement.Text = \" del mes n\u00famero \"\n$find.Forward = $true\n$find.Wrap = 0\n\n# wdReplaceOne (1): only the first (left-most) match inside this paragraph's\n# range gets replaced, which is the \" de\" right before \"{{ fecha_mes }}\".\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1) | Out-Null\n"}
